$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.508.34'
$ws.Range('E2').Value = '  -2.23%  '
$ws.Range('D3').Value = '3.703.23'
$ws.Range('E3').Value = '  -2.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '693.44'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.39'
$ws.Range('E6').Value = '  -5.13%  '
$ws.Range('D7').Value = '3.702.02'
$ws.Range('E7').Value = '  -2.96%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  -4.39%  '
$ws.Range('E10').Value = '  -7.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.38'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('E13').Value = '  -4.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.58'
$ws.Range('E14').Value = '  -6.67%  '
$ws.Range('D15').Value = '4.325.64'
$ws.Range('E15').Value = '  -3.01%  '
$ws.Range('D16').Value = '3.697.71'
$ws.Range('E16').Value = '  -3.72%  '
$ws.Range('D17').Value = '69.574.93'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.35'
$ws.Range('E19').Value = '  -6.85%  '
$ws.Range('E20').Value = '  -7.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '482.34'
$ws.Range('E21').Value = '  -5.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.02'
$ws.Range('E22').Value = '  -6.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.668'
$ws.Range('E23').Value = '  -7.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.14'
$ws.Range('E24').Value = '  -4.71%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000132'
$ws.Range('E25').Value = '  -8.83%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.846.20'
$ws.Range('E26').Value = '  -3.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.46'
$ws.Range('E28').Value = '  -4.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.57'
$ws.Range('E29').Value = '  -8.23%  '
$ws.Range('E30').Value = '  -9.52%  '
$ws.Range('E31').Value = '  -9.43%  '
$ws.Range('E32').Value = '  -7.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.09'
$ws.Range('E33').Value = '  -6.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.14'
$ws.Range('E34').Value = '  -6.67%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.168'
$ws.Range('E35').Value = '  -3.22%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').Value = '3.667.68'
$ws.Range('E37').Value = '  -2.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.53'
$ws.Range('E38').Value = '  -6.96%  '
$ws.Range('E39').Value = '  +6.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.33'
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0937'
$ws.Range('E41').Value = '  -7.42%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  -6.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '163.43'
$ws.Range('E45').Value = '  -4.46%  '
$ws.Range('E46').Value = '  -2.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '30.25'
$ws.Range('E47').Value = '  +3.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.84'
$ws.Range('E48').Value = '  -13.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.17'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('E51').Value = '  -7.91%  '
